# Worklog.xlsx update: add day-2 crawler-spider entry, goal/lines-of-code
# columns and a small "+ / -" legend row, per the commit:
#   "Main and Spider modification" -- program now extracts title+url from
#   the company page, then follows each article url to fetch date/time and
#   full text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$xlVCenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$xlShiftUp = [Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp

# ---------------------------------------------------------------------
# 1. Cell text / numeric content (order matters: it drives the order new
#    strings are appended to sharedStrings.xml, which the target diff
#    fixes at indices 3..9, so we add them in that exact sequence).
# ---------------------------------------------------------------------

# Row 4 - day 2 (2019-06-20) work-done text -> sharedStrings[3]
$ws.Range("C4").Value = "Basic spider built that crawls https://economictimes.indiatimes.com . It looks for news title, url, time/date and article in the company specific section of the website. It then saves this datato a csv file after cleaning the string a little bit."

# Row 1 headers -> sharedStrings[4], [5]
$ws.Range("G1").Value = "Goal for tomorrow"
$ws.Range("D1").Value = "Lines of Code"

# Row 3 - day 1 (2019-06-19) -> sharedStrings[6], [7]
$ws.Range("C3").Value = "Studied the process of crawling data from any webpage, watched python tutorials, built a basic web-scraper by following a youtube tutorial."
$ws.Range("G3").Value = "To Build a basic web scraper that scrapes data of one company from one website and start working from there."

# Row 2 legend -> sharedStrings[8], [9]
$ws.Range("D2").Value = "+"
$ws.Range("E2").Value = "-"

# Numeric / date values (not shared strings). The old row-2 entry
# (A2=1, B2=43636) moves down to row 4 as day 2; row 2 becomes the
# "+ / -" legend row and row 3 is the new day-1 entry.
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 43635
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 43636
